$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CargoData")
$ws.Activate()

for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value = $cell.Value2 * 100
}

$ws.Range("E20").Select()
